$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone row 11 structure (values + formatting + blank placeholder cells) into row 12
$ws.Range("A11:AX11").Copy($ws.Range("A12"))

# Cells that had data in row 11 but must be blank in row 12
$ws.Range("F12").NumberFormat = "General"
$ws.Range("F12").ClearContents()
$ws.Range("AB12").NumberFormat = "General"
$ws.Range("AB12").ClearContents()
$ws.Range("AT12").NumberFormat = "General"
$ws.Range("AT12").ClearContents()

# Overwrite with the 2021 data
$ws.Range("A12").Value = "2021年"
$ws.Range("B12").Value = 2
$ws.Range("C12").Value = 24
$ws.Range("D12").Value = 86
$ws.Range("E12").Value = 4945
$ws.Range("G12").Value = 221
$ws.Range("I12").Value = 68
$ws.Range("J12").Value = 5
$ws.Range("L12").Value = 4361
$ws.Range("O12").Value = 39
$ws.Range("P12").Value = 12
$ws.Range("Q12").Value = 2
$ws.Range("U12").Value = 10
$ws.Range("Y12").Value = 939
$ws.Range("AA12").Value = 11668
$ws.Range("AC12").Value = 34
$ws.Range("AE12").Value = 77
$ws.Range("AF12").Value = 116
$ws.Range("AI12").Value = 21
$ws.Range("AK12").Value = 404
$ws.Range("AN12").Value = 7
$ws.Range("AO12").Value = 7
$ws.Range("AP12").Value = 56
$ws.Range("AQ12").Value = 11
$ws.Range("AS12").Value = 10
$ws.Range("AV12").Value = 151
$ws.Range("AX12").Value = 60
